# Update the division-practice answer table.
# The sheet contains one table; data rows are 1, 5, 9, 13, 17 (5 cells each),
# the rows in between are blank "work space" rows that stay untouched.
# Several source strings repeat (e.g. "43÷7=6, 1" and "32÷8=4, 0" each occur
# twice), so the cells are addressed directly by (row, column) instead of via
# Find/Replace to avoid ambiguity.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "44÷4=11, 0"
$t.Cell(1, 2).Range.Text = "19÷5=3, 4"
$t.Cell(1, 3).Range.Text = "77÷8=9, 5"
$t.Cell(1, 4).Range.Text = "28÷6=4, 4"
$t.Cell(1, 5).Range.Text = "68÷9=7, 5"

$t.Cell(5, 1).Range.Text = "68÷2=34, 0"
$t.Cell(5, 2).Range.Text = "84÷5=16, 4"
$t.Cell(5, 3).Range.Text = "62÷9=6, 8"
$t.Cell(5, 4).Range.Text = "84÷2=42, 0"
$t.Cell(5, 5).Range.Text = "75÷4=18, 3"

$t.Cell(9, 1).Range.Text = "78÷8=9, 6"
$t.Cell(9, 2).Range.Text = "24÷9=2, 6"
$t.Cell(9, 3).Range.Text = "27÷3=9, 0"
$t.Cell(9, 4).Range.Text = "48÷3=16, 0"
$t.Cell(9, 5).Range.Text = "66÷5=13, 1"

$t.Cell(13, 1).Range.Text = "22÷3=7, 1"
$t.Cell(13, 2).Range.Text = "47÷8=5, 7"
$t.Cell(13, 3).Range.Text = "73÷4=18, 1"
$t.Cell(13, 4).Range.Text = "26÷5=5, 1"
$t.Cell(13, 5).Range.Text = "77÷7=11, 0"

$t.Cell(17, 1).Range.Text = "79÷5=15, 4"
$t.Cell(17, 2).Range.Text = "59÷2=29, 1"
$t.Cell(17, 3).Range.Text = "74÷7=10, 4"
$t.Cell(17, 4).Range.Text = "65÷9=7, 2"
$t.Cell(17, 5).Range.Text = "76÷2=38, 0"

Write-Host "Replacements applied."
